$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    if ($val -match '^-?\d+(\.\d+)?$') {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $val
}

Set-TextValue "D2" "64.435.36"
Set-TextValue "E2" "  +0.23%  "
Set-TextValue "D3" "3.514.02"
Set-TextValue "E3" "  +0.25%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "592.19"
Set-TextValue "E5" "  +1.47%  "
Set-TextValue "D6" "134.53"
Set-TextValue "E6" "  -0.35%  "
Set-TextValue "D8" "0.489"
Set-TextValue "E8" "  +0.20%  "
Set-TextValue "D9" "7.57"
Set-TextValue "E9" "  +6.54%  "
Set-TextValue "D10" "0.125"
Set-TextValue "E10" "  +0.20%  "
Set-TextValue "D11" "0.387"
Set-TextValue "E11" "  +3.71%  "
Set-TextValue "D12" "4.109.90"
Set-TextValue "E12" "  +0.20%  "
Set-TextValue "E13" "  +1.59%  "
Set-TextValue "E14" "  +1.08%  "
Set-TextValue "D15" "3.513.81"
Set-TextValue "E15" "  +0.20%  "
Set-TextValue "D16" "25.91"
Set-TextValue "E16" "  -1.48%  "
Set-TextValue "D17" "64.419.21"
Set-TextValue "E17" "  +0.18%  "
Set-TextValue "D18" "9.94"
Set-TextValue "E18" "  +2.01%  "
Set-TextValue "E19" "  +3.37%  "
Set-TextValue "D20" "13.63"
Set-TextValue "E20" "  -1.36%  "
Set-TextValue "D21" "394.66"
Set-TextValue "E21" "  +2.92%  "
Set-TextValue "D22" "0.576"
Set-TextValue "E22" "  +1.60%  "
Set-TextValue "D23" "3.654.01"
Set-TextValue "E23" "  +0.24%  "
Set-TextValue "E24" "  +1.12%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.03%  "
Set-TextValue "E26" "  +0.05%  "
Set-TextValue "E27" "  +3.04%  "
Set-TextValue "D28" "0.998"
Set-TextValue "E28" "  -0.16%  "
Set-TextValue "D29" "7.42"
Set-TextValue "E29" "  -2.21%  "
Set-TextValue "E30" "  +2.27%  "
Set-TextValue "E31" "  +0.36%  "
Set-TextValue "E32" "  -6.27%  "
Set-TextValue "E33" "  +7.12%  "
Set-TextValue "D34" "3.541.01"
Set-TextValue "E34" "  +0.51%  "
Set-TextValue "D36" "23.41"
Set-TextValue "E36" "  -0.54%  "
Set-TextValue "D37" "5.36"
Set-TextValue "E38" "  +1.85%  "
Set-TextValue "E39" "  +1.20%  "
Set-TextValue "D40" "167.26"
Set-TextValue "E40" "  +1.84%  "
Set-TextValue "D41" "0.0793"
Set-TextValue "E41" "  +1.54%  "
Set-TextValue "E42" "  +0.28%  "
Set-TextValue "E43" "  +0.02%  "
Set-TextValue "D44" "25.46"
Set-TextValue "E44" "  -1.54%  "
Set-TextValue "D45" "4.46"
Set-TextValue "E45" "  +1.33%  "
Set-TextValue "E46" "  +3.80%  "
Set-TextValue "E47" "  -2.94%  "
Set-TextValue "E48" "  +0.75%  "
Set-TextValue "D49" "2.398.28"
Set-TextValue "E49" "  -3.02%  "
Set-TextValue "D50" "0.899"
Set-TextValue "E50" "  -2.34%  "
Set-TextValue "E51" "  +0.03%  "
